$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Fix typo in header: "Paxkage size" -> "Package size"
$ws.Range("E1").Value = "Package size"

# Update the selected/active cell to E1 (reflecting the review of the header edit)
$ws.Range("E1").Select()
